$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1166.125
$ws.Range("I111").Value = 824
$ws.Range("J111").Value = 1508.25
$ws.Range("K111").Value = 2472
$ws.Range("L111").Value = 4524.75
$ws.Range("M111").Value = 595
$ws.Range("N111").Value = -10658.75
$ws.Range("H112").Value = 23810712
$ws.Range("I112").Value = 333333800
$ws.Range("J112").Value = 1244.6923
$ws.Range("K112").Value = 1000001400
$ws.Range("L112").Value = 3734.0769
$ws.Range("M112").Value = -1000000292
$ws.Range("N112").Value = -5950.0769
$ws.Range("H125").Value = 2961.2
$ws.Range("I125").Value = 1232
$ws.Range("J125").Value = 5555
$ws.Range("K125").Value = 11088
$ws.Range("L125").Value = 49995
$ws.Range("M125").Value = -8628
$ws.Range("N125").Value = -54915
$ws.Range("H129").Value = 828.12
$ws.Range("I129").Value = 316.58334
$ws.Range("J129").Value = 897.875
$ws.Range("K129").Value = 949.7500200000001
$ws.Range("L129").Value = 2693.625
$ws.Range("M129").Value = 4050.24998
$ws.Range("N129").Value = -12693.625
$ws.Range("H137").Value = 1362728.4
$ws.Range("I137").Value = 2647164.2
$ws.Range("J137").Value = 2737.353
$ws.Range("K137").Value = 7941492.600000001
$ws.Range("L137").Value = 8212.059000000001
$ws.Range("M137").Value = -7938942.600000001
$ws.Range("N137").Value = -13312.059
$ws.Range("H138").Value = 5582.59
$ws.Range("I138").Value = 1058.1666
$ws.Range("J138").Value = 6575.756
$ws.Range("K138").Value = 3174.4998
$ws.Range("L138").Value = 19727.268
$ws.Range("M138").Value = 1965.5002
$ws.Range("N138").Value = -30007.268

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9668
$ws.Range("I3").Value = 1005
$ws.Range("J3").Value = 13999.5
$ws.Range("K3").Value = 1005
$ws.Range("L3").Value = 13999.5
$ws.Range("M3").Value = -890
$ws.Range("N3").Value = -14229.5
$ws.Range("H61").Value = 3259.1738
$ws.Range("I61").Value = 1039.3529
$ws.Range("J61").Value = 9548.666999999999
$ws.Range("K61").Value = 1039.3529
$ws.Range("L61").Value = 9548.666999999999
$ws.Range("M61").Value = -827.3529000000001
$ws.Range("N61").Value = -9972.666999999999
$ws.Range("H74").Value = 5083.8076
$ws.Range("I74").Value = 7739
$ws.Range("K74").Value = 7739
$ws.Range("M74").Value = -6865
$ws.Range("H77").Value = 5083.8076
$ws.Range("I77").Value = 7739
$ws.Range("K77").Value = 38695
$ws.Range("M77").Value = -34327
$ws.Range("H92").Value = 26137.5
$ws.Range("J92").Value = 26137.5
$ws.Range("L92").Value = 26137.5
$ws.Range("N92").Value = -31129.5
$ws.Range("H132").Value = 1423.2192
$ws.Range("I132").Value = 884.78723
$ws.Range("J132").Value = 2396.5386
$ws.Range("K132").Value = 2654.36169
$ws.Range("L132").Value = 7189.6158
$ws.Range("M132").Value = -124.3616900000002
$ws.Range("N132").Value = -12249.6158
$ws.Range("H136").Value = 3259.1738
$ws.Range("I136").Value = 1039.3529
$ws.Range("J136").Value = 9548.666999999999
$ws.Range("K136").Value = 3118.0587
$ws.Range("L136").Value = 28646.001
$ws.Range("M136").Value = -568.0587000000005
$ws.Range("N136").Value = -33746.001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4863.184
$ws.Range("I134").Value = 1554.1482
$ws.Range("K134").Value = 4662.444600000001
$ws.Range("M134").Value = -2127.444600000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3013.8518
$ws.Range("I31").Value = 1190
$ws.Range("J31").Value = 5293.6665
$ws.Range("K31").Value = 1190
$ws.Range("L31").Value = 5293.6665
$ws.Range("M31").Value = -895
$ws.Range("N31").Value = -5883.6665
$ws.Range("H34").Value = 3013.8518
$ws.Range("I34").Value = 1190
$ws.Range("J34").Value = 5293.6665
$ws.Range("K34").Value = 1190
$ws.Range("L34").Value = 5293.6665
$ws.Range("M34").Value = -988
$ws.Range("N34").Value = -5697.6665
$ws.Range("H58").Value = 2431.3171
$ws.Range("J58").Value = 6821.357
$ws.Range("L58").Value = 6821.357
$ws.Range("N58").Value = -7227.357
$ws.Range("H132").Value = 2931.3684
$ws.Range("I132").Value = 2466.4849
$ws.Range("K132").Value = 7399.4547
$ws.Range("M132").Value = -4869.4547
$ws.Range("H134").Value = 2317.2632
$ws.Range("I134").Value = 1240.6154
$ws.Range("K134").Value = 3721.8462
$ws.Range("M134").Value = -1186.8462
$ws.Range("H136").Value = 2431.3171
$ws.Range("J136").Value = 6821.357
$ws.Range("L136").Value = 20464.071
$ws.Range("N136").Value = -25564.071

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 3000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 3000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H97").Value = 1194.2
$ws.Range("I97").Value = 1080
$ws.Range("J97").Value = 1270.3334
$ws.Range("K97").Value = 1080
$ws.Range("L97").Value = 1270.3334
$ws.Range("M97").Value = -584
$ws.Range("N97").Value = -2262.3334
$ws.Range("H132").Value = 2744.682
$ws.Range("I132").Value = 1628.5
$ws.Range("J132").Value = 3382.5
$ws.Range("K132").Value = 4885.5
$ws.Range("L132").Value = 10147.5
$ws.Range("M132").Value = -2355.5
$ws.Range("N132").Value = -15207.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 28560
$ws.Range("J64").Value = 28560
$ws.Range("L64").Value = 28560
$ws.Range("N64").Value = -29010
$ws.Range("H67").Value = 28560
$ws.Range("J67").Value = 28560
$ws.Range("L67").Value = 28560
$ws.Range("N67").Value = -30120
$ws.Range("H94").Value = 31818.625
$ws.Range("J94").Value = 31818.625
$ws.Range("L94").Value = 31818.625
$ws.Range("N94").Value = -33170.625
$ws.Range("H132").Value = 22330
$ws.Range("I132").Value = 28534.666
$ws.Range("K132").Value = 85603.99800000001
$ws.Range("M132").Value = -83073.99800000001
$ws.Range("H136").Value = 3989.4285
$ws.Range("I136").Value = 2062.8333
$ws.Range("J136").Value = 5434.375
$ws.Range("K136").Value = 6188.499899999999
$ws.Range("L136").Value = 16303.125
$ws.Range("M136").Value = -3638.499899999999
$ws.Range("N136").Value = -21403.125

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 17530
$ws.Range("J98").Value = 17530
$ws.Range("L98").Value = 17530
$ws.Range("N98").Value = -23520
$ws.Range("H122").Value = 2570.1
$ws.Range("I122").Value = 1612.08
$ws.Range("K122").Value = 4836.24
$ws.Range("M122").Value = -2386.24
$ws.Range("H132").Value = 1684.439
$ws.Range("I132").Value = 1117.3939
$ws.Range("K132").Value = 3352.1817
$ws.Range("M132").Value = -822.1817000000001
$ws.Range("H136").Value = 2967.5527
$ws.Range("I136").Value = 2043.1482
$ws.Range("K136").Value = 6129.444600000001
$ws.Range("M136").Value = -3579.444600000001
